$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.738.91"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.625.17"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").Value = "'0.992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.81%  "
$ws.Range("D5").Value = "'210.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").Value = "'23.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").Value = "'0.257"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.854.47"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "1.624.14"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "'0.560"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "'65.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "27.734.48"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "'231.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "'7.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'10.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.40%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("E24").Value = "  -4.42%  "
$ws.Range("D25").Value = "'153.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "'6.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").Value = "'0.110"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").Value = "'15.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "'0.993"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "'3.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").Value = "'3.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "1.394.04"
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("D35").Value = "'1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.97%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").Value = "'0.0169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "'0.556"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "'0.864"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").Value = "'1.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("D42").Value = "'0.992"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'66.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "'5.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "'2.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "1.765.23"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").Value = "'87.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("D50").Value = "'0.0994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "'0.0505"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
